$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update KISW score for row 2 (Steve Coup)
$ws.Range("E2").Value = 78
# Update recalculated TOTAL and AVERAGE for row 2
$ws.Range("K2").Value = 545
$ws.Range("L2").Value = 68.125

# Update KISW score for row 3 (Fredrick Ndote)
$ws.Range("E3").Value = 60
# Update recalculated TOTAL and AVERAGE for row 3
$ws.Range("K3").Value = 541
$ws.Range("L3").Value = 67.625
